# Switch example group numbers
# Replace the example team member names (Alice, Bob, Claire, David, Elaine)
# with the real group's names, and add a 6th member (Martin) in a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$names = @{
    12 = "Veselin"
    13 = "Rawda"
    14 = "Hannah"
    15 = "Mirit"
    16 = "Bogdana"
}

foreach ($row in $names.Keys) {
    $cell = $ws.Range("B$row")
    $cell.Value = $names[$row]
    $cell.Font.Color = 0
}

# New row for the 6th group member, keeping the same rating/marker pattern
$ws.Range("B17").Value = "Martin"
$ws.Range("B17").Font.Color = 0
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 1

# Update the active selection like in the target workbook
$ws.Range("F15").Select()
